# SMS Scripts dry run
# Refresh the "dry run" SMS sample data (phone numbers, dates, CurrentTime
# stamps) across all 4 worksheets of the WF6_Send_SMS_ETR_M_Test workbook.
#
# Notes:
#  - Several of the new values are purely numeric / date-shaped strings
#    (e.g. "2512700123", "2024-05-25", "0364350912"). A plain
#    Range.Value = "..." assignment would get auto-coerced by Excel into a
#    number/date, which both loses data (leading zeros) and swaps the
#    cell's style out from under us. Prefixing the literal with a leading
#    apostrophe forces text entry, matching how Excel's "Text" cells work.
#  - Forcing text entry like that nudges the cell onto the workbook's
#    built-in quote-prefix style. To keep each cell's original style
#    (s="2" for all the cells touched here) we snapshot an untouched,
#    already-style-"2" cell (AQ2, on every sheet) with Copy(), write the
#    new value into the target, then PasteSpecial formats-only so the
#    target's look is restored without touching its (now text) content.
#    (Copying the target cell onto itself is not reliable once its value
#    has already been overwritten, so a separate donor cell is used.)

function Set-TextValue {
    param($Worksheet, $Address, $Value, $StyleDonorAddress)

    $target = $Worksheet.Range($Address)
    $donor = $Worksheet.Range($StyleDonorAddress)

    # Snapshot the desired style before we touch the target's contents.
    $donor.Copy()

    # Leading apostrophe forces this to be stored as text even though the
    # value looks like a number or a date.
    $target.Value = "'" + $Value

    # xlPasteFormats restores the original style (number format, font,
    # alignment, quote-prefix flag, ...) that Value= just perturbed.
    $target.PasteSpecial(-4122)
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

Set-TextValue $ws1 "F2"  "2512700123" "AQ2"
Set-TextValue $ws1 "N2"  "2024-05-25" "AQ2"
$ws1.Range("P2").Value = "2024-04-28 05:00:00 PM"
Set-TextValue $ws1 "AC2" "2024-05-25" "AQ2"
Set-TextValue $ws1 "AE2" "4004784379" "AQ2"
Set-TextValue $ws1 "AK2" "3"          "AQ2"
Set-TextValue $ws1 "AT2" "0364350912" "AQ2"
Set-TextValue $ws1 "AX2" "1172955431" "AQ2"
$ws1.Range("AZ2").Value = "CT: Sat, May 25, 2024 at 2:20 PM"

# ---------------------------------------------------------------------
# Sheet2, Sheet3, Sheet4 share the same set of edits, except the final
# "CurrentTime" stamp written into AZ2 differs per sheet.
# ---------------------------------------------------------------------
$az2ValueBySheet = @{ 2 = "CT: Sat, May 25, 2024 at 2:28 PM"; 3 = "CT: Sat, May 25, 2024 at 2:37 PM"; 4 = "CT: Sat, May 25, 2024 at 2:37 PM" }

foreach ($sheetIndex in 2, 3, 4) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    $ws.Range("AZ1").Value = "CurrentTime"

    Set-TextValue $ws "F2"  "2512700123" "AQ2"
    Set-TextValue $ws "N2"  "2024-05-25" "AQ2"
    $ws.Range("O2").Value = "02:35:55 PM"
    $ws.Range("P2").Value = "2024-04-28 05:00:00 PM"
    Set-TextValue $ws "AC2" "2024-05-25" "AQ2"
    Set-TextValue $ws "AE2" "4004784379" "AQ2"
    Set-TextValue $ws "AT2" "0364350912" "AQ2"
    Set-TextValue $ws "AX2" "1172955431" "AQ2"

    # AZ2 changes style (s="1" -> s="2") in addition to its value. Its new
    # text ("CT: ...") is not numeric-looking, so no apostrophe/Value-then-
    # restore dance is needed for the *value* -- but pasting the format
    # first (instead of going through Set-TextValue's quote-prefix step)
    # avoids ever materialising a transient "wrapText + quotePrefix" style
    # combination that would otherwise leak into styles.xml as unused.
    $ws.Range("AQ2").Copy()
    $ws.Range("AZ2").PasteSpecial(-4122)
    $ws.Range("AZ2").Value = $az2ValueBySheet[$sheetIndex]
}
